$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "79.732.10"
$ws.Range("E2").Value = "  +4.36%  "
$ws.Range("D3").Value = "3.203.30"
$ws.Range("E3").Value = "  +5.05%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'205.67"
$ws.Range("E5").Value = "  +2.51%  "
$ws.Range("D6").Value = "'636.67"
$ws.Range("E6").Value = "  +1.96%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.246"
$ws.Range("E8").Value = "  +19.63%  "
$ws.Range("D9").Value = "'0.605"
$ws.Range("E9").Value = "  +10.32%  "
$ws.Range("D10").Value = "3.200.68"
$ws.Range("E10").Value = "  +5.11%  "
$ws.Range("D11").Value = "'0.622"
$ws.Range("E11").Value = "  +40.75%  "
$ws.Range("E12").Value = "  +32.19%  "
$ws.Range("E13").Value = "  +3.42%  "
$ws.Range("E14").Value = "  +3.74%  "
$ws.Range("D15").Value = "3.793.43"
$ws.Range("E15").Value = "  +5.53%  "
$ws.Range("D16").Value = "'32.41"
$ws.Range("E16").Value = "  +11.20%  "
$ws.Range("D17").Value = "79.635.90"
$ws.Range("E17").Value = "  +4.35%  "
$ws.Range("D18").Value = "3.200.24"
$ws.Range("E18").Value = "  +5.10%  "
$ws.Range("D19").Value = "'14.68"
$ws.Range("E19").Value = "  +8.38%  "
$ws.Range("D20").Value = "'9.40"
$ws.Range("E20").Value = "  +4.20%  "
$ws.Range("D21").Value = "'447.39"
$ws.Range("E21").Value = "  +19.20%  "
$ws.Range("E22").Value = "  +27.82%  "
$ws.Range("D23").Value = "'5.29"
$ws.Range("E23").Value = "  +20.87%  "
$ws.Range("D24").Value = "'4.82"
$ws.Range("E24").Value = "  +10.13%  "
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").Value = "3.370.43"
$ws.Range("E25").Value = "  +5.63%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "'77.85"
$ws.Range("E26").Value = "  +6.22%  "
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").Value = "'10.89"
$ws.Range("E27").Value = "  +11.16%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "'0.0000122"
$ws.Range("E29").Value = "  +11.34%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'9.23"
$ws.Range("E30").Value = "  +11.21%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.49"
$ws.Range("E32").Value = "  +6.11%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'536.23"
$ws.Range("E33").Value = "  +8.58%  "
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").Value = "'2.03"
$ws.Range("E34").Value = "  +4.18%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.148"
$ws.Range("E35").Value = "  +27.34%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "'23.39"
$ws.Range("E36").Value = "  +13.08%  "
$ws.Range("B37").Value = "Cronos"
$ws.Range("C37").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D37").Value = "'0.124"
$ws.Range("E37").Value = "  +18.69%  "
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").Value = "'0.411"
$ws.Range("E39").Value = "  +7.08%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "'164.61"
$ws.Range("E40").Value = "  +1.11%  "
$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D41").Value = "'20.05"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'192.56"
$ws.Range("E42").Value = "  +0.89%  "
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D44").Value = "'5.56"
$ws.Range("E44").Value = "  +8.08%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'1.83"
$ws.Range("E45").Value = "  +10.78%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.806"
$ws.Range("E46").Value = "  -0.73%  "
$ws.Range("D47").Value = "'2.66"
$ws.Range("E47").Value = "  +7.93%  "
$ws.Range("B48").Value = "ImmutableX"
$ws.Range("C48").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D48").Value = "'1.33"
$ws.Range("E48").Value = "  +4.43%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "'43.70"
$ws.Range("E49").Value = "  +3.90%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'25.78"
$ws.Range("E50").Value = "  +15.40%  "
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").Value = "'0.640"
$ws.Range("E51").Value = "  +5.43%  "
